$d = $word.ActiveDocument

# Locate the "User Class" attributes line ("Attributes: Username, Password")
# and insert "Name, " right before "Username" so it reads
# "Attributes: Name, Username, Password".
$found = $d.Content.Find.Execute("Username, Password", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Name, Username, Password", 2)
